$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row to append: 06-12-2025 gold price update.
$newDate = "06-12-2025"
$newText = "The price of gold in India today is ₹13,015 per gram for 24 karat gold, ₹11,930 per gram for 22 karat gold and ₹9,761 per gram for 18 karat gold (also called 999 gold)."

# Write the date-like text through a text formula in a scratch cell, then
# paste-special *values only* into the target cell. Going through a text
# formula (and a values-only paste) keeps the string a literal shared
# string instead of letting Excel auto-convert the "dd-mm-yyyy" looking
# text into a date serial number (which would also fabricate a brand new
# cell style). The scratch cell is fully cleared afterwards so it leaves
# no trace in the saved sheet (dimension, styles, etc. stay untouched).
$scratch = $ws.Range("Z1")
$scratch.Formula = '="' + $newDate + '"'
$scratch.Copy()
$ws.Range("A81").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$scratch.Clear()

# Column B is plain text already, so a direct value assignment is safe.
$ws.Range("B81").Value = $newText
